$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section boundaries so that the edits below
# only touch bullets inside that section (some bullet text is duplicated verbatim
# elsewhere in the document, e.g. under PROFESSIONAL EXPERIENCE).
$sectionStart = 0
$sectionEnd = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $sectionStart = $i
    }
    elseif ($sectionStart -gt 0 -and $sectionEnd -eq 0 -and $t -like "*TECHNICAL SKILLS*") {
        $sectionEnd = $i
    }
}

function Set-ParagraphText($doc, $rangeStart, $rangeEnd, $oldText, $newText) {
    for ($i = $rangeStart; $i -le $rangeEnd; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($oldText + "`r")) {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

function Remove-ParagraphInRange($doc, $rangeStart, $rangeEnd, $targetText) {
    for ($i = $rangeEnd; $i -ge $rangeStart; $i--) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($targetText + "`r")) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# 1) Rewrite the three bullets that are changed in place.
$r1 = Set-ParagraphText $d $sectionStart $sectionEnd `
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations" `
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"

$r2 = Set-ParagraphText $d $sectionStart $sectionEnd `
    "• Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets" `
    "• `$4.7M savings enabled nonprofit access"

$r3 = Set-ParagraphText $d $sectionStart $sectionEnd `
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" `
    "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

# 2) Remove the two bullets that are deleted entirely.
$r4 = Remove-ParagraphInRange $d $sectionStart $sectionEnd `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"

$r5 = Remove-ParagraphInRange $d $sectionStart $sectionEnd `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"

# 3) Rewrite the final bullet (recompute the section end since deletions shift indices).
$sectionEnd = 0
$count = $d.Paragraphs.Count
for ($i = $sectionStart; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*TECHNICAL SKILLS*") {
        $sectionEnd = $i
        break
    }
}

$r6 = Set-ParagraphText $d $sectionStart $sectionEnd `
    "• Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy" `
    "• Real-time collaboration at national scale"

if (-not ($r1 -and $r2 -and $r3 -and $r4 -and $r5 -and $r6)) {
    Write-Output "WARNING: one or more Key Achievements edits did not find their target paragraph (r1=$r1 r2=$r2 r3=$r3 r4=$r4 r5=$r5 r6=$r6)"
}
